# Split the "So far I see ... methods" paragraph (which originally had the
# run text broken mid-word around a bookmark) into four paragraphs matching
# the author's new journal entries, preserving the list formatting on the
# first paragraph and relocating the _GoBack bookmark and the page-break run.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("So far I see the class itself can be abstract")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the target paragraph ('So far I see the class itself can be abstract...')"
}

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00267465" w:rsidRDefault="004D1261" w:rsidP="004D1261"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>So far I see the class itself can be abstract, as well as the methods</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Since I clearly cannot have virtual and override fields, I’ll have to reconsider. </w:t></w:r><w:r><w:t xml:space="preserve"> I thought I could use the magic of OO to avoid writing that line in every file, but it looks like no. At least not for tonight. Hmm… But I do for sure want the functionality of a static class variable… and I still imagine there will be benefits to this inheritance relationship I’ve created, and I don’t want to rename that variable differently for each class. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>so</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> now I must consider how to avoid </w:t></w:r><w:r w:rsidR="00267465"><w:br w:type="page"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# InsertXML replaces the contents of the target Range with the supplied
# WordprocessingML, so this single call turns the one paragraph into the
# four paragraphs described by the new content.
$null = $target.Range.InsertXML($xml)
